$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 29591
$ws.Range("E2").Value = 1595
$ws.Range("F2").Value = 1595
$ws.Range("G2").Value = 1292
$ws.Range("H2").Value = 931
$ws.Range("I2").Value = 364
$ws.Range("J2").Value = 566
$ws.Range("K2").Value = 21801
$ws.Range("L2").Value = 12112
$ws.Range("M2").Value = 9690
$ws.Range("N2").Value = 4580
$ws.Range("O2").Value = 5110
$ws.Range("P2").Value = 371
$ws.Range("Q2").Value = 1832
$ws.Range("R2").Value = -1701
$ws.Range("S2").Value = -447
$ws.Range("T2").Value = 1105
$ws.Range("U2").Value = 728
$ws.Range("V2").Value = 7048
$ws.Range("W2").Value = 5.39
$ws.Range("X2").Value = 3.15
$ws.Range("Y2").Value = 8.19
$ws.Range("Z2").Value = 4.36
$ws.Range("AA2").Value = 125
$ws.Range("AB2").Value = 1195.12
$ws.Range("AC2").Value = 981
$ws.Range("AD2").Value = 17.12
$ws.Range("AE2").Value = 12335
$ws.Range("AF2").Value = 1.36
$ws.Range("AG2").Value = 170
$ws.Range("AH2").Value = 1.01
$ws.Range("AI2").Value = 17.35
$ws.Range("AJ2").Value = 36212538
$ws.Range("D3").Value = 30517
$ws.Range("E3").Value = 1330
$ws.Range("F3").Value = 1330
$ws.Range("G3").Value = 929
$ws.Range("H3").Value = 530
$ws.Range("I3").Value = 244
$ws.Range("J3").Value = 286
$ws.Range("K3").Value = 24804
$ws.Range("L3").Value = 14824
$ws.Range("M3").Value = 9980
$ws.Range("N3").Value = 4719
$ws.Range("O3").Value = 5261
$ws.Range("P3").Value = 371
$ws.Range("Q3").Value = 1165
$ws.Range("R3").Value = -2251
$ws.Range("S3").Value = 1861
$ws.Range("T3").Value = 1896
$ws.Range("U3").Value = -732
$ws.Range("V3").Value = 9081
$ws.Range("W3").Value = 4.36
$ws.Range("X3").Value = 1.74
$ws.Range("Y3").Value = 5.25
$ws.Range("Z3").Value = 2.27
$ws.Range("AA3").Value = 148.53
$ws.Range("AB3").Value = 1231.26
$ws.Range("AC3").Value = 658
$ws.Range("AD3").Value = 29.11
$ws.Range("AE3").Value = 12712
$ws.Range("AF3").Value = 1.51
$ws.Range("AG3").Value = 170
$ws.Range("AH3").Value = 0.89
$ws.Range("AI3").Value = 25.88
$ws.Range("AJ3").Value = 36212538
$ws.Range("D4").Value = 33181
$ws.Range("E4").Value = 1361
$ws.Range("F4").Value = 1361
$ws.Range("G4").Value = 989
$ws.Range("H4").Value = 672
$ws.Range("I4").Value = 287
$ws.Range("J4").Value = 385
$ws.Range("K4").Value = 26893
$ws.Range("L4").Value = 16391
$ws.Range("M4").Value = 10502
$ws.Range("N4").Value = 4963
$ws.Range("O4").Value = 5539
$ws.Range("P4").Value = 371
$ws.Range("Q4").Value = 782
$ws.Range("R4").Value = -2150
$ws.Range("S4").Value = 903
$ws.Range("T4").Value = 1206
$ws.Range("U4").Value = -424
$ws.Range("V4").Value = 10383
$ws.Range("W4").Value = 4.1
$ws.Range("X4").Value = 2.03
$ws.Range("Y4").Value = 5.93
$ws.Range("Z4").Value = 2.6
$ws.Range("AA4").Value = 156.07
$ws.Range("AB4").Value = 1288.31
$ws.Range("AC4").Value = 774
$ws.Range("AD4").Value = 13.31
$ws.Range("AE4").Value = 13368
$ws.Range("AF4").Value = 0.77
$ws.Range("AG4").Value = 170
$ws.Range("AH4").Value = 1.65
$ws.Range("AI4").Value = 22.01
$ws.Range("AJ4").Value = 36212538
$ws.Range("D5").Value = 33837
$ws.Range("E5").Value = 1145
$ws.Range("F5").Value = 1145
$ws.Range("G5").Value = 990
$ws.Range("H5").Value = 594
$ws.Range("I5").Value = 260
$ws.Range("J5").Value = 334
$ws.Range("K5").Value = 26511
$ws.Range("L5").Value = 15779
$ws.Range("M5").Value = 10733
$ws.Range("N5").Value = 5088
$ws.Range("O5").Value = 5645
$ws.Range("P5").Value = 371
$ws.Range("Q5").Value = 1853
$ws.Range("R5").Value = -835
$ws.Range("S5").Value = -1213
$ws.Range("T5").Value = 1428
$ws.Range("U5").Value = 425
$ws.Range("V5").Value = 9335
$ws.Range("W5").Value = 3.38
$ws.Range("X5").Value = 1.75
$ws.Range("Y5").Value = 5.16
$ws.Range("Z5").Value = 2.22
$ws.Range("AA5").Value = 147.02
$ws.Range("AB5").Value = 1347.88
$ws.Range("AC5").Value = 699
$ws.Range("AD5").Value = 14.23
$ws.Range("AE5").Value = 13705
$ws.Range("AF5").Value = 0.73
$ws.Range("AG5").Value = 180
$ws.Range("AH5").Value = 1.81
$ws.Range("AI5").Value = 25.78
$ws.Range("AJ5").Value = 36212538
$ws.Range("D6").Value = 33980
$ws.Range("E6").Value = 1329
$ws.Range("F6").Value = 1329
$ws.Range("G6").Value = 976
$ws.Range("H6").Value = 456
$ws.Range("I6").Value = 235
$ws.Range("K6").Value = 26296
$ws.Range("L6").Value = 15348
$ws.Range("M6").Value = 10948
$ws.Range("N6").Value = 5227
$ws.Range("P6").Value = 371
$ws.Range("Q6").Value = 1427
$ws.Range("R6").Value = -996
$ws.Range("S6").Value = -245
$ws.Range("T6").Value = 1556
$ws.Range("U6").Value = -130
$ws.Range("V6").Value = 9283
$ws.Range("W6").Value = 3.91
$ws.Range("X6").Value = 1.34
$ws.Range("Y6").Value = 4.55
$ws.Range("Z6").Value = 1.73
$ws.Range("AA6").Value = 140.19
$ws.Range("AB6").Value = 1397.5
$ws.Range("AC6").Value = 632
$ws.Range("AD6").Value = 12.44
$ws.Range("AE6").Value = 14078
$ws.Range("AF6").Value = 0.56
$ws.Range("AG6").Value = 190
$ws.Range("AH6").Value = 2.42
$ws.Range("AI6").Value = 30.11
$ws.Range("AJ6").Value = 36212538

# Estimate rows 7-9 (2020E, 2021E, 2022E) no longer have financial data -
# clear all columns D through AJ, leaving only the A/B/C identifier columns.
$ws.Range("D7:AJ9").ClearContents()

Write-Host "Done updating Daesang Holdings IFRS data"
